# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1!A1 : update the "Conversión del día" text block ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$oldText = $ws1.Range("A1").Text
$newText = $oldText -replace "1000 Bs = 7\.07 = 28500\.71 pesos", "1000 Bs = 7.11 = 28492.18 pesos"
$newText = $newText -replace "28500\.71 pesos = 7\.09 = 974\.93 Bs", "28492.18 pesos = 7.08 = 954.91 Bs"
$ws1.Range("A1").Value = $newText

# --- tasas sheet : update N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 140.6
$ws2.Range("O10").Value = 4006
$ws2.Range("N12").Value = 4026
$ws2.Range("O12").Value = 134.93
